# Updates Leve profit-tracking figures (currentAveragePrice / NQ / HQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H-N)
# on several rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR
# sheets, reflecting refreshed market-board pricing data pulled in by
# the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 960.8333
$ws.Range("J28").Value = 400
$ws.Range("L28").Value = 400
$ws.Range("N28").Value = -1370

# Row 40
$ws.Range("H40").Value = 5529.6294
$ws.Range("I40").Value = 1016.6667
$ws.Range("K40").Value = 1016.6667
$ws.Range("M40").Value = -841.6667

# Row 62
$ws.Range("H62").Value = 2667.6667
$ws.Range("I62").Value = 3001.5
$ws.Range("K62").Value = 3001.5
$ws.Range("M62").Value = -2377.5

# Row 65
$ws.Range("H65").Value = 2667.6667
$ws.Range("I65").Value = 3001.5
$ws.Range("K65").Value = 15007.5
$ws.Range("M65").Value = -11887.5

# Row 80
$ws.Range("H80").Value = 899.2857
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502

# Row 83
$ws.Range("H83").Value = 899.2857
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 492

# Row 86
$ws.Range("H86").Value = 2633.6667
$ws.Range("I86").Value = 916
$ws.Range("K86").Value = 916
$ws.Range("M86").Value = 207

# Row 89
$ws.Range("H89").Value = 2633.6667
$ws.Range("I89").Value = 916
$ws.Range("K89").Value = 4580
$ws.Range("M89").Value = 1036

# Row 106
$ws.Range("H106").Value = 3996.5
$ws.Range("I106").Value = 3996.5
$ws.Range("K106").Value = 3996.5
$ws.Range("M106").Value = -3365.5

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1914.5
$ws.Range("I86").Value = 1986
$ws.Range("K86").Value = 1986
$ws.Range("M86").Value = -863

# Row 89
$ws.Range("H89").Value = 1914.5
$ws.Range("I89").Value = 1986
$ws.Range("K89").Value = 9930
$ws.Range("M89").Value = -4314

# Row 94
$ws.Range("H94").Value = 1880.5
$ws.Range("I94").Value = 1880.5
$ws.Range("K94").Value = 1880.5
$ws.Range("M94").Value = -1429.5

# Row 107
$ws.Range("H107").Value = 3814
$ws.Range("I107").Value = 3814
$ws.Range("K107").Value = 3814
$ws.Range("M107").Value = -1894

$ws = $wb.Worksheets.Item("CRP")
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 16
$ws.Range("H16").Value = 675
$ws.Range("I16").Value = 675
$ws.Range("K16").Value = 675
$ws.Range("M16").Value = -388

# Row 31
$ws.Range("H31").Value = 4033.3635
$ws.Range("J31").Value = 4699.4
$ws.Range("L31").Value = 4699.4
$ws.Range("N31").Value = -5289.4

# Row 34
$ws.Range("H34").Value = 4033.3635
$ws.Range("J34").Value = 4699.4
$ws.Range("L34").Value = 4699.4
$ws.Range("N34").Value = -5103.4

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# Row 51
$ws.Range("H51").Value = 19300
$ws.Range("J51").Value = 19300
$ws.Range("L51").Value = 19300
$ws.Range("N51").Value = -20772

# Row 61
$ws.Range("H61").Value = 19300
$ws.Range("J61").Value = 19300
$ws.Range("L61").Value = 19300
$ws.Range("N61").Value = -19996

# Row 99
$ws.Range("H99").Value = 2252.75
$ws.Range("I99").Value = 12
$ws.Range("J99").Value = 2999.6667
$ws.Range("K99").Value = 12
$ws.Range("L99").Value = 2999.6667
$ws.Range("M99").Value = 1486
$ws.Range("N99").Value = -5995.6667

# Row 107
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 800
$ws.Range("L107").Value = 350
$ws.Range("M107").Value = 1120
$ws.Range("N107").Value = -4190

# Row 113
$ws.Range("H113").Value = 675
$ws.Range("I113").Value = 675
$ws.Range("K113").Value = 675
$ws.Range("M113").Value = 1495

# Row 126
$ws.Range("H126").Value = 2252.75
$ws.Range("I126").Value = 12
$ws.Range("J126").Value = 2999.6667
$ws.Range("K126").Value = 36
$ws.Range("L126").Value = 8999.000100000001
$ws.Range("M126").Value = 2434
$ws.Range("N126").Value = -13939.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2169.1428
$ws.Range("I97").Value = 1530.6666
$ws.Range("K97").Value = 1530.6666
$ws.Range("M97").Value = -1034.6666

# Row 122
$ws.Range("H122").Value = 1189.8
$ws.Range("I122").Value = 1189.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3569.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1119.4
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 9000
$ws.Range("I40").Value = 9000
$ws.Range("K40").Value = 9000
$ws.Range("M40").Value = -8864

# Row 46
$ws.Range("H46").Value = 4801.6
$ws.Range("J46").Value = 4801.6
$ws.Range("L46").Value = 4801.6
$ws.Range("N46").Value = -5177.6

# Row 55
$ws.Range("H55").Value = 2954.5454
$ws.Range("I55").Value = 300.66666
$ws.Range("J55").Value = 3949.75
$ws.Range("K55").Value = 300.66666
$ws.Range("L55").Value = 3949.75
$ws.Range("M55").Value = -127.66666
$ws.Range("N55").Value = -4295.75

# Row 68
$ws.Range("H68").Value = 1990
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 1990
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 122
$ws.Range("H122").Value = 4511
$ws.Range("I122").Value = 4663.2
$ws.Range("K122").Value = 13989.6
$ws.Range("M122").Value = -11539.6

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 24400.4
$ws.Range("J18").Value = 24999
$ws.Range("L18").Value = 24999
$ws.Range("N18").Value = -25345

# Row 96
$ws.Range("H96").Value = 8000
$ws.Range("I96").Value = 8000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 8000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -6627
$ws.Range("N96").ClearContents()

# Row 112
$ws.Range("H112").Value = 24998
$ws.Range("J112").Value = 24998
$ws.Range("L112").Value = 24998
$ws.Range("N112").Value = -27952

# Row 122
$ws.Range("H122").Value = 5542.8887
$ws.Range("I122").Value = 5140.857
$ws.Range("J122").Value = 6950
$ws.Range("K122").Value = 15422.571
$ws.Range("L122").Value = 20850
$ws.Range("M122").Value = -12972.571
$ws.Range("N122").Value = -25750

# Row 126
$ws.Range("H126").Value = 1736.875
$ws.Range("J126").Value = 1281
$ws.Range("L126").Value = 3843
$ws.Range("N126").Value = -8783
